$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FM")
$ws.Select()

$ws.Range("D2").Value = 15
$ws.Range("E2").Value = 42
$ws.Range("F2").Value = 280

$ws.Range("D3").Value = 16
$ws.Range("E3").Value = 25
$ws.Range("F3").Value = 280

$ws.Range("D4").Value = 12
$ws.Range("E4").Value = 49
$ws.Range("F4").Value = 280

$ws.Range("D4").Select()
